# Add a new "test_xlr_n_percent" column (K) to the "table_test_" table,
# populated with "n (pct%)" style labels for rows 3..34, plus an empty
# but identically-styled cell for the trailing blank row 35 - mirroring
# how the existing text columns (test_r_char / test_r_factor / etc.) look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Header for the new column
$ws.Range("K2").Value = "test_xlr_n_percent"

# Data rows 3..34 -> "1 (3%)" ... "32 (100%)"
$ws.Range("K3").Value  = "1 (3%)"
$ws.Range("K4").Value  = "2 (6%)"
$ws.Range("K5").Value  = "3 (9%)"
$ws.Range("K6").Value  = "4 (12%)"
$ws.Range("K7").Value  = "5 (16%)"
$ws.Range("K8").Value  = "6 (19%)"
$ws.Range("K9").Value  = "7 (22%)"
$ws.Range("K10").Value = "8 (25%)"
$ws.Range("K11").Value = "9 (28%)"
$ws.Range("K12").Value = "10 (31%)"
$ws.Range("K13").Value = "11 (34%)"
$ws.Range("K14").Value = "12 (38%)"
$ws.Range("K15").Value = "13 (41%)"
$ws.Range("K16").Value = "14 (44%)"
$ws.Range("K17").Value = "15 (47%)"
$ws.Range("K18").Value = "16 (50%)"
$ws.Range("K19").Value = "17 (53%)"
$ws.Range("K20").Value = "18 (56%)"
$ws.Range("K21").Value = "19 (59%)"
$ws.Range("K22").Value = "20 (62%)"
$ws.Range("K23").Value = "21 (66%)"
$ws.Range("K24").Value = "22 (69%)"
$ws.Range("K25").Value = "23 (72%)"
$ws.Range("K26").Value = "24 (75%)"
$ws.Range("K27").Value = "25 (78%)"
$ws.Range("K28").Value = "26 (81%)"
$ws.Range("K29").Value = "27 (84%)"
$ws.Range("K30").Value = "28 (88%)"
$ws.Range("K31").Value = "29 (91%)"
$ws.Range("K32").Value = "30 (94%)"
$ws.Range("K33").Value = "31 (97%)"
$ws.Range("K34").Value = "32 (100%)"

# Style the first data cell to match the other text columns in the table
# (calibri font, right/bottom aligned, General number format), then copy
# that formatting down across the rest of the new column, including the
# trailing blank row (35), so only a single new style entry is created.
$seed = $ws.Range("K3")
$seed.HorizontalAlignment = -4152
$seed.VerticalAlignment = -4107
$seed.Font.Name = "calibri"

$seed.Copy()
$ws.Range("K4:K35").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Grow the table so it now spans A2:K34 and includes the new column.
$tbl.Resize($ws.Range("A2:K34"))
$ws.Range("K2").Value = "test_xlr_n_percent"
